$d = $word.ActiveDocument

# 1) Locate the paragraph "1 - une page pour le détail  d'un todo." and update
#    its text in place to "2 - Une page pour le détail  d'un todo." (a plain
#    in-place Range.Text assignment, so the existing paragraph/run keep their
#    original formatting/fidelity).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "1 - une page pour le détail*") {
        $target = $p
        break
    }
}
$target.Range.Text = "2 - Une page pour le détail  d’un todo."

# 2) The document always ends with one empty paragraph. Grow it in place into
#    three paragraphs - "3 - Un formulaire pour ajouter un todo." followed by
#    two empty ones - by assigning text containing embedded carriage returns
#    (paragraph marks) to its Range. This both adds the new
#    "3 - Un formulaire pour ajouter un todo." paragraph and restores the
#    trailing empty paragraph the document ends with, while every new run
#    inherits the same Helvetica/24-half-point formatting used throughout.
$cr = [char]13
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Text = "3 - Un formulaire pour ajouter un todo.$cr$cr"
